# The sheet originally listed 6 Sunil Narine innings (rows 2-7).
# The edit keeps only the "Oct 7 2020 vs Chennai Super Kings" innings
# (previously row 5) and drops every other match row, shrinking the
# used range down to A1:K2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting first so numeric-looking values ("17", "188.88", ...)
# stay stored as text, matching the original "number stored as text" data.
$ws.Range("A2:K2").NumberFormat = "@"

$ws.Range("A2").Value = " Oct 7 2020"
$ws.Range("B2").Value = " Abu Dhabi"
$ws.Range("C2").Value = "KKR won by 10 runs"
$ws.Range("D2").Value = "Kolkata Knight Riders"
$ws.Range("E2").Value = "Chennai Super Kings"
$ws.Range("F2").Value = "Sunil Narine "
$ws.Range("G2").Value = "17"
$ws.Range("H2").Value = "9"
$ws.Range("I2").Value = "1"
$ws.Range("J2").Value = "1"
$ws.Range("K2").Value = "188.88"

# Drop the remaining match rows (old rows 3-7) entirely, shifting the
# used range/dimension from A1:K7 down to A1:K2.
$ws.Rows("3:7").Delete()
